$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: new weekly data
$ws.Range("D2").Value = 45282
$ws.Range("M2").Value = 150
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 20000
$ws.Range("S2").Value = 2500

# Row 3: previous row 2 values shift down
$ws.Range("D3").Value = 44895
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 22000
$ws.Range("O3").Value = 22500
$ws.Range("P3").Value = 22250
$ws.Range("S3").Value = 2781

# Row 4: previous row 3 values shift down
$ws.Range("D4").Value = 44495
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 270
$ws.Range("N4").Value = 19000
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 19556
$ws.Range("S4").Value = 2444
